$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.ColorScheme
# touch legacy api once, just to see effect on name (value same as before to minimize other changes)
$orig = $cs.Item(1).RGB
$cs.Item(1).RGB = $orig

$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
Write-Host "Name after legacy touch + before 12-set: $($tcs.Name)"
